$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.7432489451418292
$ws.Range("C2").Value = 0.3275137255068472
$ws.Range("E2").Value = 0.1432144595564671
$ws.Range("F2").Value = 0.4443680307746121
$ws.Range("G2").Value = 0.3768493239248585
$ws.Range("H2").Value = 0.5425939907787836
$ws.Range("L2").Value = 0.1873566682000529
$ws.Range("M2").Value = 0.1773456444137409
$ws.Range("N2").Value = 1.231579395578748
$ws.Range("O2").Value = 1.780361369930148
$ws.Range("B3").Value = 0.66954623271792
$ws.Range("C3").Value = 0.3223261300460081
$ws.Range("E3").Value = 0.1446139863887028
$ws.Range("F3").Value = 0.3878228170618172
$ws.Range("G3").Value = 0.3777371580639723
$ws.Range("H3").Value = 0.5466047476233982
$ws.Range("L3").Value = 0.1846137862165946
$ws.Range("M3").Value = 0.1651921032376222
$ws.Range("N3").Value = 1.237017196794682
$ws.Range("O3").Value = 1.790392869546423
$ws.Range("B4").Value = 0.6243102249437129
$ws.Range("C4").Value = 0.319150910934539
$ws.Range("E4").Value = 0.1455307761656144
$ws.Range("F4").Value = 0.3531389305169483
$ws.Range("G4").Value = 0.3785928258561739
$ws.Range("H4").Value = 0.54933331371371
$ws.Range("L4").Value = 0.1830184604745781
$ws.Range("M4").Value = 0.1577708680731433
$ws.Range("N4").Value = 1.240815561463641
$ws.Range("O4").Value = 1.797758625050875
$ws.Range("B5").Value = 0.6058820113603929
$ws.Range("C5").Value = 0.3178596524363684
$ws.Range("E5").Value = 0.1459188416525059
$ws.Range("F5").Value = 0.3390132514313251
$ws.Range("G5").Value = 0.3790194973410479
$ws.Range("H5").Value = 0.5505121278861367
$ws.Range("L5").Value = 0.1823907598339503
$ws.Range("M5").Value = 0.1547572144467004
$ws.Range("N5").Value = 1.242479205144647
$ws.Range("O5").Value = 1.801063408439802
$ws.Range("B6").Value = 0.6028224125718964
$ws.Range("C6").Value = 0.317645405809543
$ws.Range("E6").Value = 0.1459841537581834
$ws.Range("F6").Value = 0.336668177824194
$ws.Range("G6").Value = 0.3790950524847716
$ws.Range("H6").Value = 0.5507119105623133
$ws.Range("L6").Value = 0.1822878861687585
$ws.Range("M6").Value = 0.154257444065518
$ws.Range("N6").Value = 1.242762451403813
$ws.Range("O6").Value = 1.801630471306439
$ws.Range("B7").Value = 0.6240616702117165
$ws.Range("C7").Value = 0.3191334855537349
$ws.Range("E7").Value = 0.1455359511608227
$ws.Range("F7").Value = 0.3529483938344953
$ws.Range("G7").Value = 0.3785982645005816
$ws.Range("H7").Value = 0.5493489406970724
$ws.Range("L7").Value = 0.1830099042577302
$ws.Range("M7").Value = 0.1577301818272545
$ws.Range("N7").Value = 1.240837528856126
$ws.Range("O7").Value = 1.797801967190324
$ws.Range("B8").Value = 0.7178333585451639
$ws.Range("C8").Value = 0.3257230631415098
$ws.Range("E8").Value = 0.1436850958959868
$ws.Range("F8").Value = 0.4248636149813478
$ws.Range("G8").Value = 0.3770909329118126
$ws.Range("H8").Value = 0.5439217189654642
$ws.Range("L8").Value = 0.1863925322630919
$ws.Range("M8").Value = 0.1731467083183134
$ws.Range("N8").Value = 1.233359106966525
$ws.Range("O8").Value = 1.783569740737505
$ws.Range("B9").Value = 0.9018093312985229
$ws.Range("C9").Value = 0.3387181407347555
$ws.Range("E9").Value = 0.1405109702784477
$ws.Range("F9").Value = 0.5661985755041457
$ws.Range("G9").Value = 0.3766043254576772
$ws.Range("H9").Value = 0.5353880623315916
$ws.Range("L9").Value = 0.1937278303496015
$ws.Range("M9").Value = 0.2036960296856947
$ws.Range("N9").Value = 1.222330794884257
$ws.Range("O9").Value = 1.765241773877776
$ws.Range("B10").Value = 1.036976067903652
$ws.Range("C10").Value = 0.3483025880940716
$ws.Range("E10").Value = 0.1384556492740403
$ws.Range("F10").Value = 0.6702781546542269
$ws.Range("G10").Value = 0.3777602487548819
$ws.Range("H10").Value = 0.5304029709303251
$ws.Range("L10").Value = 0.1995421333308798
$ws.Range("M10").Value = 0.2263251390662049
$ws.Range("N10").Value = 1.216433734424157
$ws.Range("O10").Value = 1.757632363173883
$ws.Range("B11").Value = 1.098455904810066
$ws.Range("C11").Value = 0.3526693012730533
$ws.Range("E11").Value = 0.1375805206145282
$ws.Range("F11").Value = 0.7176906081379002
$ws.Range("G11").Value = 0.3786165109947746
$ws.Range("H11").Value = 0.528413831723725
$ws.Range("L11").Value = 0.2022789438646839
$ws.Range("M11").Value = 0.2366579888118849
$ws.Range("N11").Value = 1.214227550079087
$ws.Range("O11").Value = 1.755445617550947
$ws.Range("B12").Value = 1.121734326974376
$ws.Range("C12").Value = 0.3543236775471001
$ws.Range("E12").Value = 0.1372577242603459
$ws.Range("F12").Value = 0.7356546913071611
$ws.Range("G12").Value = 0.3789883977532895
$ws.Range("H12").Value = 0.5277006399441433
$ws.Range("L12").Value = 0.2033284470411587
$ws.Range("M12").Value = 0.2405761411733991
$ws.Range("N12").Value = 1.213460438048799
$ws.Range("O12").Value = 1.754801106247243
$ws.Range("B13").Value = 1.116721044008898
$ws.Range("C13").Value = 0.3539673452262946
$ws.Range("E13").Value = 0.1373268620870016
$ws.Range("F13").Value = 0.7317853510981394
$ws.Range("G13").Value = 0.3789061845789519
$ws.Range("H13").Value = 0.5278524574012238
$ws.Range("L13").Value = 0.2031018349010623
$ws.Range("M13").Value = 0.239732063986061
$ws.Range("N13").Value = 1.213622613631301
$ws.Range("O13").Value = 1.75493174545673
$ws.Range("B14").Value = 1.100371096349306
$ws.Range("C14").Value = 0.3528053927890085
$ws.Range("E14").Value = 0.137553791747635
$ws.Range("F14").Value = 0.7191683204515869
$ws.Range("G14").Value = 0.3786461508844781
$ws.Range("H14").Value = 0.5283543543661011
$ws.Range("L14").Value = 0.2023650243662445
$ws.Range("M14").Value = 0.2369802322264221
$ws.Range("N14").Value = 1.214163071087285
$ws.Range("O14").Value = 1.755388912972904
$ws.Range("B15").Value = 1.090355891292802
$ws.Range("C15").Value = 0.3520937621954232
$ws.Range("E15").Value = 0.1376939118556297
$ws.Range("F15").Value = 0.7114413442032514
$ws.Range("G15").Value = 0.3784930803724791
$ws.Range("H15").Value = 0.5286669963879405
$ws.Range("L15").Value = 0.2019154146039881
$ws.Range("M15").Value = 0.2352953421615496
$ws.Range("N15").Value = 1.214503008976962
$ws.Range("O15").Value = 1.755692853666773
$ws.Range("B16").Value = 1.032957924609207
$ws.Range("C16").Value = 0.3480173336748322
$ws.Range("E16").Value = 0.1385140444063557
$ws.Range("F16").Value = 0.6671810134426437
$ws.Range("G16").Value = 0.3777109496560769
$ws.Range("H16").Value = 0.5305385711529311
$ws.Range("L16").Value = 0.1993651181719969
$ws.Range("M16").Value = 0.2256506218042631
$ws.Range("N16").Value = 1.216587484501318
$ws.Range("O16").Value = 1.757800945677729
$ws.Range("B17").Value = 0.9977429798786943
$ws.Range("C17").Value = 0.3455181740333586
$ws.Range("E17").Value = 0.1390324902808553
$ws.Range("F17").Value = 0.6400460337125793
$ws.Range("G17").Value = 0.3773158566094565
$ws.Range("H17").Value = 0.5317580710946856
$ws.Range("L17").Value = 0.1978240659945101
$ws.Range("M17").Value = 0.2197436488098035
$ws.Range("N17").Value = 1.217988131077391
$ws.Range("O17").Value = 1.759420877572268
$ws.Range("B18").Value = 0.9774875925822357
$ws.Range("C18").Value = 0.3440813640669091
$ws.Range("E18").Value = 0.1393363203355609
$ws.Range("F18").Value = 0.6244449056556647
$ws.Range("G18").Value = 0.3771197035684821
$ws.Range("H18").Value = 0.5324857193706407
$ws.Range("L18").Value = 0.1969463439293406
$ws.Range("M18").Value = 0.2163497766440869
$ws.Range("N18").Value = 1.218838599547709
$ws.Range("O18").Value = 1.760472596561272
$ws.Range("B19").Value = 0.9706293943953597
$ws.Range("C19").Value = 0.3435950000059762
$ws.Range("E19").Value = 0.139440159872394
$ws.Range("F19").Value = 0.619163680173358
$ws.Range("G19").Value = 0.3770586261533282
$ws.Range("H19").Value = 0.5327365929299077
$ws.Range("L19").Value = 0.1966506503513727
$ws.Range("M19").Value = 0.2152013063867741
$ws.Range("N19").Value = 1.219134262947222
$ws.Range("O19").Value = 1.760849289188911
$ws.Range("B20").Value = 1.001491752223046
$ws.Range("C20").Value = 0.3457841487562945
$ws.Range("E20").Value = 0.1389767178388377
$ws.Range("F20").Value = 0.642933953830422
$ws.Range("G20").Value = 0.3773546959220795
$ws.Range("H20").Value = 0.5316255392194051
$ws.Range("L20").Value = 0.197987218898561
$ws.Range("M20").Value = 0.2203720789028125
$ws.Range("N20").Value = 1.217834389039112
$ws.Range("O20").Value = 1.759236014100281
$ws.Range("B21").Value = 1.105173556223974
$ws.Range("C21").Value = 0.3531466660875395
$ws.Range("E21").Value = 0.1374869037960469
$ws.Range("F21").Value = 0.7228739723491628
$ws.Range("G21").Value = 0.3787212351996629
$ws.Range("H21").Value = 0.5282058481460012
$ws.Range("L21").Value = 0.2025810876787517
$ws.Range("M21").Value = 0.2377883693254716
$ws.Range("N21").Value = 1.214002472854489
$ws.Range("O21").Value = 1.755249648037989
$ws.Range("B22").Value = 1.172919431616663
$ws.Range("C22").Value = 0.3579630640243181
$ws.Range("E22").Value = 0.1365633228877865
$ws.Range("F22").Value = 0.7751780083420101
$ws.Range("G22").Value = 0.379892079673894
$ws.Range("H22").Value = 0.5262043313898488
$ws.Range("L22").Value = 0.2056599568335997
$ws.Range("M22").Value = 0.24920185408061
$ws.Range("N22").Value = 1.211896243870768
$ws.Range("O22").Value = 1.753714377650482
$ws.Range("B23").Value = 1.136764151722161
$ws.Range("C23").Value = 0.3553920988039323
$ws.Range("E23").Value = 0.1370516746870658
$ws.Range("F23").Value = 0.7472568307830727
$ws.Range("G23").Value = 0.3792417248798756
$ws.Range("H23").Value = 0.5272512215384921
$ws.Range("L23").Value = 0.2040097309374289
$ws.Range("M23").Value = 0.2431075135547047
$ws.Range("N23").Value = 1.212984007734818
$ws.Range("O23").Value = 1.754435790463617
$ws.Range("B24").Value = 0.9997969628035435
$ws.Range("C24").Value = 0.3456639016022649
$ws.Range("E24").Value = 0.139001914590601
$ws.Range("F24").Value = 0.6416283278902171
$ws.Range("G24").Value = 0.3773370401434022
$ws.Range("H24").Value = 0.5316853742044003
$ws.Range("L24").Value = 0.1979134317695639
$ws.Range("M24").Value = 0.2200879590190254
$ws.Range("N24").Value = 1.217903754949248
$ws.Range("O24").Value = 1.759319215862718
$ws.Range("B25").Value = 0.8520350443448592
$ws.Range("C25").Value = 0.3351955973779326
$ws.Range("E25").Value = 0.141320996577039
$ws.Range("F25").Value = 0.5279251897347166
$ws.Range("G25").Value = 0.3764707718717233
$ws.Range("H25").Value = 0.5374709704719649
$ws.Range("L25").Value = 0.1916685842924153
$ws.Range("M25").Value = 0.1953985702206822
$ws.Range("N25").Value = 1.224926118756144
$ws.Range("O25").Value = 1.769172653912634
